# Update the header of the results table (this is the first column header
# of Table2, "MNIST Dataset Experiment"). Editing the cell's value also
# keeps the table's ListColumn name in sync with the header cell, since the
# table header row is backed by the worksheet cell itself.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "MNIST Dataset Experiment (96 Samples)"

# Clear the stale selection (previously parked on J21, well outside the
# used range) by moving the active selection back to the top-left cell.
$ws.Range("A1").Select()
